$d = $word.ActiveDocument

# --- First paragraph: pPr changes (pBdr + ind) ---
$p1 = $d.Paragraphs.Item(1)

# Add a paragraph border (top/left/bottom/right, 5pt space) to the first paragraph.
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# --- First paragraph: run/text changes ---
# Replace the placeholder id text, then drop the trailing run (a lone space).
$markerOld = "**ID__AFFARS_pgi_5315_topic_14__ID**"
$markerNew = "**ID__AFFARS_AFICC_PGI_5315_407_90__ID**"

$rng = $d.Range($p1.Range.Start, $p1.Range.Start + $markerOld.Length)
$rng.Text = $markerNew

# Remove the remaining trailing-space run so the paragraph ends right after the marker.
$tailStart = $p1.Range.Start + $markerNew.Length
$tailEnd = $p1.Range.End - 1
if ($tailEnd -gt $tailStart) {
    $tailRng = $d.Range($tailStart, $tailEnd)
    $tailRng.Delete()
}
